$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # "Transactions" sheet

# ---------------------------------------------------------------------------
# 1. Shift the old last data row (row 14 - the Disbursement row) down to row
#    15, preserving its formatting, then populate row 14 with a new
#    transaction row (renumbered sequence continues: ... 308, 307, 306).
# ---------------------------------------------------------------------------
$ws.Range("A14:J14").Copy()
$ws.Range("A15:J15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(15, 1).Value = 306
$ws.Cells.Item(15, 2).Value = "Head Office"
$ws.Cells.Item(15, 3).Value = 41917
$ws.Cells.Item(15, 4).Value = "Disbursement"
$ws.Cells.Item(15, 5).Value = 5000
$ws.Cells.Item(15, 10).Value = 5000

# New row 14 (same look/style as the row above it, row 13)
$ws.Range("A13:J13").Copy()
$ws.Range("A14:J14").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(14, 1).Value = 307
$ws.Cells.Item(14, 2).Value = "Head Office"
$ws.Cells.Item(14, 3).Value = 41943
$ws.Cells.Item(14, 4).Value = "Accrual"
$ws.Cells.Item(14, 5).Value = 44.38
$ws.Cells.Item(14, 10).Value = 0

# ---------------------------------------------------------------------------
# 2. Correct the transaction IDs / dates / types / amounts / running balances
#    for the existing rows 2-13.
# ---------------------------------------------------------------------------

# Row 2 (Income Posting)
$ws.Cells.Item(2, 1).Value = 318
$ws.Cells.Item(2, 10).Value = 5300.19

# Row 3 (Accrual)
$ws.Cells.Item(3, 1).Value = 317

# Row 4 (Income Posting)
$ws.Cells.Item(4, 1).Value = 316
$ws.Cells.Item(4, 10).Value = 5246.72

# Row 5 (Accrual)
$ws.Cells.Item(5, 1).Value = 315

# Row 6 (Income Posting)
$ws.Cells.Item(6, 1).Value = 314
$ws.Cells.Item(6, 10).Value = 5198.5 + 0.36

# Row 7 (Accrual)
$ws.Cells.Item(7, 1).Value = 313

# Row 8 - becomes an Income Posting row with a running balance
$ws.Cells.Item(8, 1).Value = 332
$ws.Cells.Item(8, 4).Value = "Income Posting"
$ws.Cells.Item(8, 5).Value = 52.28
$ws.Range("J2:J2").Copy()
$ws.Range("J8:J8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(8, 10).Value = 5146.41

# Row 9 - becomes an Accrual row with zero running balance
$ws.Cells.Item(9, 1).Value = 331
$ws.Cells.Item(9, 4).Value = "Accrual"
$ws.Cells.Item(9, 5).Value = 0.36
$ws.Range("J3:J3").Copy()
$ws.Range("J9:J9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(9, 10).Value = 0

# Row 10 - becomes an Accrual row with zero running balance
$ws.Cells.Item(10, 1).Value = 311
$ws.Cells.Item(10, 3).Value = 42004
$ws.Cells.Item(10, 4).Value = "Accrual"
$ws.Cells.Item(10, 5).Value = 51.92
$ws.Range("J3:J3").Copy()
$ws.Range("J10:J10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(10, 10).Value = 0

# Row 11 - becomes an Income Posting row with a running balance
$ws.Cells.Item(11, 1).Value = 310
$ws.Cells.Item(11, 4).Value = "Income Posting"
$ws.Range("J2:J2").Copy()
$ws.Range("J11:J11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(11, 10).Value = 5094.13

# Row 12 - becomes an Accrual row with zero running balance
$ws.Cells.Item(12, 1).Value = 309
$ws.Cells.Item(12, 3).Value = 41973
$ws.Cells.Item(12, 4).Value = "Accrual"
$ws.Cells.Item(12, 5).Value = 49.75
$ws.Range("J3:J3").Copy()
$ws.Range("J12:J12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(12, 10).Value = 0

# Row 13 - becomes an Income Posting row with a running balance
$ws.Cells.Item(13, 1).Value = 308
$ws.Cells.Item(13, 4).Value = "Income Posting"
$ws.Range("J2:J2").Copy()
$ws.Range("J13:J13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(13, 10).Value = 5044.38

# ---------------------------------------------------------------------------
# 3. Make "Transactions" the active sheet (was "ChargesTab" before).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("G5").Select()
